$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "TC16_17" worksheet right after "TC_05_06", before
#    "Error" (this matches the new sheetId=6 / shifted rIds in the diff).
# ---------------------------------------------------------------------
$tc0506 = $wb.Worksheets.Item("TC_05_06")
$tc1617 = $wb.Worksheets.Add($null, $tc0506)
$tc1617.Name = "TC16_17"

# ---------------------------------------------------------------------
# 2. Populate the new sheet with its header row + two data rows.
# ---------------------------------------------------------------------
$headers = @("Code", "Coverage_Type", "Mileage_Band", "Class", "Term", "Value", "Surcharge", "Option", "Deductibles")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $tc1617.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$dataRow = @("SNE", "Powertrain", "0-60", "'2", "24/24", "Y", "Y", "N", "N")
for ($i = 0; $i -lt $dataRow.Length; $i++) {
    $tc1617.Cells.Item(2, $i + 1).Value = $dataRow[$i]
    $tc1617.Cells.Item(3, $i + 1).Value = $dataRow[$i]
}

# ---------------------------------------------------------------------
# 3. Column widths on the new sheet.
# ---------------------------------------------------------------------
$tc1617.Columns.Item(2).ColumnWidth = 13.665
$tc1617.Columns.Item(3).ColumnWidth = 14.165
$tc1617.Columns.Item(9).ColumnWidth = 10.665

# Selection on the new sheet.
$tc1617.Range("B11").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Re-activate "TC_05_06" and update its selection / scroll position.
# ---------------------------------------------------------------------
$tc0506.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$tc0506.Range("S10").Select() | Out-Null
